$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 97
$ws.Range("C2").Value = "flower/flower095.png"
$ws.Range("D2").Value = "krachen"
$ws.Range("E2").Value = "flower"
$ws.Range("B3").Value = 110
$ws.Range("C3").Value = "dog/dog117.png"
$ws.Range("D3").Value = "schmecken"
$ws.Range("E3").Value = "dog"
$ws.Range("B4").Value = 67
$ws.Range("C4").Value = "dog/dog069.png"
$ws.Range("D4").Value = "mieten"
$ws.Range("E4").Value = "dog"
$ws.Range("B5").Value = 57
$ws.Range("C5").Value = "flower/flower102.png"
$ws.Range("D5").Value = "segeln"
$ws.Range("E5").Value = "flower"
$ws.Range("B6").Value = 25
$ws.Range("C6").Value = "dog/dog098.png"
$ws.Range("D6").Value = "langen"
$ws.Range("E6").Value = "dog"
$ws.Range("B7").Value = 42
$ws.Range("C7").Value = "dog/dog079.png"
$ws.Range("D7").Value = "fliehen"
$ws.Range("E7").Value = "dog"
$ws.Range("B8").Value = 9
$ws.Range("C8").Value = "dog/dog073.png"
$ws.Range("D8").Value = "währen"
$ws.Range("E8").Value = "dog"
$ws.Range("B9").Value = 41
$ws.Range("C9").Value = "flower/flower067.png"
$ws.Range("D9").Value = "fliegen"
$ws.Range("E9").Value = "flower"
$ws.Range("B10").Value = 91
$ws.Range("C10").Value = "dog/dog112.png"
$ws.Range("D10").Value = "kaufen"
$ws.Range("E10").Value = "dog"
$ws.Range("B11").Value = 35
$ws.Range("C11").Value = "flower/flower076.png"
$ws.Range("D11").Value = "kehren"
$ws.Range("E11").Value = "flower"
$ws.Range("B12").Value = 66
$ws.Range("C12").Value = "flower/flower085.png"
$ws.Range("D12").Value = "schätzen"
$ws.Range("E12").Value = "flower"
$ws.Range("B13").Value = 104
$ws.Range("C13").Value = "flower/flower117.png"
$ws.Range("D13").Value = "hupen"
$ws.Range("E13").Value = "flower"
$ws.Range("B14").Value = 68
$ws.Range("C14").Value = "flower/flower111.png"
$ws.Range("D14").Value = "laufen"
$ws.Range("E14").Value = "flower"
$ws.Range("B15").Value = 63
$ws.Range("C15").Value = "flower/flower104.png"
$ws.Range("D15").Value = "gelten"
$ws.Range("E15").Value = "flower"
$ws.Range("B16").Value = 121
$ws.Range("C16").Value = "dog/dog120.png"
$ws.Range("D16").Value = "füllen"
$ws.Range("E16").Value = "dog"
$ws.Range("B17").Value = 32
$ws.Range("C17").Value = "flower/flower099.png"
$ws.Range("D17").Value = "opfern"
$ws.Range("E17").Value = "flower"
$ws.Range("B18").Value = 37
$ws.Range("C18").Value = "dog/dog072.png"
$ws.Range("D18").Value = "sieben"
$ws.Range("E18").Value = "dog"
$ws.Range("B19").Value = 109
$ws.Range("C19").Value = "flower/flower113.png"
$ws.Range("D19").Value = "formen"
$ws.Range("E19").Value = "flower"
$ws.Range("B20").Value = 52
$ws.Range("C20").Value = "dog/dog108.png"
$ws.Range("D20").Value = "gründen"
$ws.Range("E20").Value = "dog"
$ws.Range("B21").Value = 99
$ws.Range("C21").Value = "flower/flower121.png"
$ws.Range("D21").Value = "lehnen"
$ws.Range("E21").Value = "flower"
$ws.Range("B22").Value = 82
$ws.Range("C22").Value = "flower/flower073.png"
$ws.Range("D22").Value = "starten"
$ws.Range("E22").Value = "flower"
$ws.Range("B23").Value = 7
$ws.Range("C23").Value = "dog/dog088.png"
$ws.Range("D23").Value = "spielen"
$ws.Range("E23").Value = "dog"
$ws.Range("B24").Value = 17
$ws.Range("C24").Value = "dog/dog099.png"
$ws.Range("D24").Value = "liefern"
$ws.Range("E24").Value = "dog"
$ws.Range("B25").Value = 108
$ws.Range("C25").Value = "flower/flower094.png"
$ws.Range("D25").Value = "stärken"
$ws.Range("E25").Value = "flower"
$ws.Range("B26").Value = 27
$ws.Range("C26").Value = "dog/dog103.png"
$ws.Range("D26").Value = "töten"
$ws.Range("E26").Value = "dog"
$ws.Range("B27").Value = 117
$ws.Range("C27").Value = "flower/flower110.png"
$ws.Range("D27").Value = "loben"
$ws.Range("E27").Value = "flower"
$ws.Range("B28").Value = 83
$ws.Range("C28").Value = "dog/dog095.png"
$ws.Range("D28").Value = "hauen"
$ws.Range("E28").Value = "dog"
$ws.Range("B29").Value = 89
$ws.Range("C29").Value = "flower/flower071.png"
$ws.Range("D29").Value = "jubeln"
$ws.Range("E29").Value = "flower"
$ws.Range("B30").Value = 75
$ws.Range("C30").Value = "flower/flower096.png"
$ws.Range("D30").Value = "strahlen"
$ws.Range("E30").Value = "flower"
$ws.Range("B31").Value = 51
$ws.Range("C31").Value = "dog/dog092.png"
$ws.Range("D31").Value = "sondern"
$ws.Range("E31").Value = "dog"
$ws.Range("B32").Value = 126
$ws.Range("C32").Value = "dog/dog122.png"
$ws.Range("D32").Value = "saufen"
$ws.Range("E32").Value = "dog"
$ws.Range("B33").Value = 8
$ws.Range("C33").Value = "dog/dog089.png"
$ws.Range("D33").Value = "wenden"
$ws.Range("E33").Value = "dog"
